$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing row 66 (and below) down by one.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new record's data.
# Columns A, B, C, E, F, G, H, I, J are constant across every data row in this sheet,
# so copy them from the (now shifted) row 67 which still holds the original values.
$ws.Range("A66").Value = $ws.Range("A67").Value2
$ws.Range("B66").Value = $ws.Range("B67").Value2
$ws.Range("C66").Value = $ws.Range("C67").Value2
$ws.Range("D66").Value = 44880
$ws.Range("E66").Value = $ws.Range("E67").Value2
$ws.Range("F66").Value = $ws.Range("F67").Value2
$ws.Range("G66").Value = $ws.Range("G67").Value2
$ws.Range("H66").Value = $ws.Range("H67").Value2
$ws.Range("I66").Value = $ws.Range("I67").Value2
$ws.Range("J66").Value = $ws.Range("J67").Value2
$ws.Range("K66").Value = "Sin especificar"
$ws.Range("L66").Value = "Primera"
$ws.Range("M66").Value = 450
$ws.Range("N66").Value = 7000
$ws.Range("O66").Value = 7500
$ws.Range("P66").Value = 7250
$ws.Range("Q66").Value = "$/bandeja 4 kilos"
$ws.Range("R66").Value = "Brasil"
$ws.Range("S66").Value = 1812
$ws.Range("T66").Value = 4
